# Update "Datos actualizados a..." (paises.xlsx) - COVID country stats refresh
# Commit: "Update countries & provincias Spain"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Footer timestamp text: 16:20 -> 16:50
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 16:50"

# 2) Estados Unidos (row 4): refreshed totals
$ws.Range("B4").Value = 189711
$ws.Range("C4").Value = 1181
$ws.Range("E4").Value = 178345
$ws.Range("G4").Value = 46
$ws.Range("H4").Value = 4099

# 3) Israel (row 21): refreshed totals
$ws.Range("E21").Value = 5342
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 23

# 4) Republica Dominicana gets new, higher totals and overtakes Islandia,
#    Mexico and Panama in the ranking, so it moves up from row 48 to row 45
#    and the three countries it passes each shift down one row.
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("B45").Value = 1284
$ws.Range("C45").Value = 175
$ws.Range("D45").Value = 9
$ws.Range("E45").Value = 1218
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 57

$ws.Range("A46").Value = "Islandia"
$ws.Range("B46").Value = 1220
$ws.Range("C46").Value = 85
$ws.Range("D46").Value = 225
$ws.Range("E46").Value = 993
$ws.Range("F46").Value = 11
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 2

$ws.Range("A47").Value = "Mexico"
$ws.Range("B47").Value = 1215
$ws.Range("C47").Value = 121
$ws.Range("D47").Value = 35
$ws.Range("E47").Value = 1151
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 29

$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 1181
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 9
$ws.Range("E48").Value = 1142
$ws.Range("F48").Value = 50
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 30

# 5) Eslovaquia (row 75): refreshed totals
$ws.Range("E75").Value = 396
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 1

# 6) Trinidad y Tobago (row 120): refreshed totals
$ws.Range("E120").Value = 83
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 5
